$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that needs to move
# from 45189 (2023-09-20) to 45190 (2023-09-21) for every data row
# (rows 2 through 310).
for ($r = 2; $r -le 310; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
